$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E6").Value = 40
$ws.Range("E7").Value = 40
$ws.Range("E8").Value = 16
$ws.Range("E9").Value = 24
$ws.Range("E10").Value = 24
$ws.Range("E20").Value = 16
$ws.Range("E21").Value = 16
$ws.Range("E22").Value = 24
$ws.Range("E23").Value = 24
$ws.Range("E24").Value = 24
$ws.Range("E25").Value = 80

[void]$ws.Range("H28").Select()
